$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 30: new "Change Password" test case row ---
$ws.Range("A30").Value = "Change Password"
$ws.Range("A30").Font.Bold = $false
$ws.Range("B30").Value = "Click on change password button"
$ws.Range("C30").Value = "Trader can change password from the trader panel"

# --- Sheet view: scroll/selection moved to A28 ---
[void]$ws.Range("A28").Select()

# --- Page setup: paper size / orientation ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
